# Se realizan cambios para sanity semilla 10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared-string values that actually change content
$ws.Range("C12").Value = "3045981670"
$ws.Range("C13").Value = "3045981684"
$ws.Range("E9").Value  = "3043209868"
$ws.Range("E10").Value = "3043209819"
$ws.Range("F10").Value = "732111193278811"

# Resize column E, add a new width for column F
$ws.Columns.Item(5).ColumnWidth = 20.8333333333
$ws.Columns.Item(6).ColumnWidth = 22.1666666667

# Move the viewport / selection to match the new view state
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F12").Select()
